$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the data held in row 4 and row 5 (product rows got re-ordered) ---
# Capture current ("before") values for both rows first. Value2 is used
# everywhere (both for reads and writes) so plain strings round-trip as-is.
$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$c4 = $ws.Range("C4").Value2
$d4 = $ws.Range("D4").Value2
$h4 = $ws.Range("H4").Value2
$i4 = $ws.Range("I4").Value2
$j4 = $ws.Range("J4").Value2
$k4 = $ws.Range("K4").Value2
$l4 = $ws.Range("L4").Value2
$m4 = $ws.Range("M4").Value2
$n4 = $ws.Range("N4").Value2

$a5 = $ws.Range("A5").Value2
$b5 = $ws.Range("B5").Value2
$c5 = $ws.Range("C5").Value2
$d5 = $ws.Range("D5").Value2
$h5 = $ws.Range("H5").Value2
$i5 = $ws.Range("I5").Value2
$j5 = $ws.Range("J5").Value2
$k5 = $ws.Range("K5").Value2
$l5 = $ws.Range("L5").Value2
$m5 = $ws.Range("M5").Value2
$n5 = $ws.Range("N5").Value2

# Columns that hold plain numeric-looking text (id, price, priceContextPrice)
# must keep their "text" storage type, so force the cell format to Text
# before writing the swapped values back in (mirrors the source data, which
# is not something Excel would otherwise infer as text).
foreach ($addr in @("A4", "H4", "K4", "A5", "H5", "K5")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Write old row 5 values into row 4 ...
$ws.Range("A4").Value2 = $a5
$ws.Range("B4").Value2 = $b5
$ws.Range("C4").Value2 = $c5
$ws.Range("D4").Value2 = $d5
$ws.Range("H4").Value2 = $h5
$ws.Range("I4").Value2 = $i5
$ws.Range("J4").Value2 = $j5
$ws.Range("K4").Value2 = $k5
$ws.Range("L4").Value2 = $l5
$ws.Range("M4").Value2 = $m5
$ws.Range("N4").Value2 = $n5

# ... and old row 4 values into row 5
$ws.Range("A5").Value2 = $a4
$ws.Range("B5").Value2 = $b4
$ws.Range("C5").Value2 = $c4
$ws.Range("D5").Value2 = $d4
$ws.Range("H5").Value2 = $h4
$ws.Range("I5").Value2 = $i4
$ws.Range("J5").Value2 = $j4
$ws.Range("K5").Value2 = $k4
$ws.Range("L5").Value2 = $l4
$ws.Range("M5").Value2 = $m4
$ws.Range("N5").Value2 = $n4

# --- Refresh the scrape timestamp (column O) for every data row ---
$newTimestamp = "2022-08-01 21:00:22"
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value2 = $newTimestamp
}
